$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing row 98: B98 should be a numeric value (3), not text
$ws.Cells.Item(98, 2).Value = 3

# Add new row 99 with annotation data for Ying Tang
$ws.Cells.Item(99, 1).Value = "Ying Tang"

# B99 stays a text "3" (matches the original data export quirk)
$ws.Cells.Item(99, 2).NumberFormat = "@"
$ws.Cells.Item(99, 2).Value = "3"

$ws.Cells.Item(99, 3).Value = "无"
$ws.Cells.Item(99, 4).Value = "DIS"
$ws.Cells.Item(99, 5).Value = "MET"
$ws.Cells.Item(99, 6).Value = "de0ccbe8-5f95-482f-b825-f58d1806a8aa"
$ws.Cells.Item(99, 7).Value = "rJTutzbA-_annotated.xlsx"
$ws.Cells.Item(99, 8).Value = "In the noiseless case, as you mention, the iterates of SGD converge linearly to the minimizer."
